# Auto-generated edit script for cryptos.xlsx update
# Applies the new Price/Volume(1h) figures and the row-shift of coins
# (LEO row removed, every following coin shifted up one row, THORChain appended).
#
# Writing these values with a plain Range.Value assignment would let Excel's
# COM layer "helpfully" reinterpret numeric-looking text (e.g. "1.00", "0.512")
# as real numbers, which would silently strip the formatting the source data
# relies on (trailing zeros, thousand-separator dots, etc). To keep every
# updated cell as literal text (matching the workbook's inlineStr cells), we
# temporarily force the cell to Text number-format before assigning the
# value, then restore the cell style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    $range = $Sheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

Set-CellText $ws "D2" "42.929.44"
Set-CellText $ws "E2" "  -1.14%  "
Set-CellText $ws "D3" "2.338.64"
Set-CellText $ws "E3" "  +1.24%  "
Set-CellText $ws "E4" "  -0.01%  "
Set-CellText $ws "D5" "306.72"
Set-CellText $ws "E5" "  -1.38%  "
Set-CellText $ws "D6" "100.32"
Set-CellText $ws "E6" "  -1.53%  "
Set-CellText $ws "E7" "  -4.38%  "
Set-CellText $ws "E8" "  +0.00%  "
Set-CellText $ws "D9" "0.512"
Set-CellText $ws "E9" "  -3.26%  "
Set-CellText $ws "D10" "34.95"
Set-CellText $ws "E10" "  -2.25%  "
Set-CellText $ws "E11" "  +0.58%  "
Set-CellText $ws "E12" "  -1.55%  "
Set-CellText $ws "E13" "  -0.10%  "
Set-CellText $ws "D14" "6.80"
Set-CellText $ws "E14" "  -2.73%  "
Set-CellText $ws "D15" "16.01"
Set-CellText $ws "E15" "  +7.05%  "
Set-CellText $ws "D16" "2.312.11"
Set-CellText $ws "E16" "  -0.09%  "
Set-CellText $ws "D17" "0.806"
Set-CellText $ws "E17" "  -0.34%  "
Set-CellText $ws "D18" "42.840.60"
Set-CellText $ws "E18" "  -1.11%  "
Set-CellText $ws "E19" "  +0.76%  "
Set-CellText $ws "D20" "0.0₃0912"
Set-CellText $ws "E20" "  -1.83%  "
Set-CellText $ws "D21" "11.71"
Set-CellText $ws "E21" "  -4.91%  "
Set-CellText $ws "E22" "  -0.15%  "
Set-CellText $ws "D23" "236.71"
Set-CellText $ws "E23" "  -1.97%  "
Set-CellText $ws "E24" "  +0.48%  "
Set-CellText $ws "E25" "  -2.08%  "
Set-CellText $ws "D26" "1.00"
Set-CellText $ws "E26" "  -0.11%  "
Set-CellText $ws "D27" "25.62"
Set-CellText $ws "E27" "  +4.17%  "
Set-CellText $ws "B28" "Toncoin"
Set-CellText $ws "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D28" "2.32"
Set-CellText $ws "E28" "  +9.51%  "
Set-CellText $ws "B29" "InjectiveProtocol"
Set-CellText $ws "C29" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-CellText $ws "D29" "35.03"
Set-CellText $ws "E29" "  -4.46%  "
Set-CellText $ws "B30" "Cosmos"
Set-CellText $ws "C30" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-CellText $ws "D30" "9.33"
Set-CellText $ws "E30" "  -3.31%  "
Set-CellText $ws "B31" "Monero"
Set-CellText $ws "C31" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText $ws "D31" "160.97"
Set-CellText $ws "E31" "  -3.80%  "
Set-CellText $ws "B32" "FirstDigitalUSD"
Set-CellText $ws "C32" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText $ws "D32" "1.00"
Set-CellText $ws "E32" "  -0.03%  "
Set-CellText $ws "B33" "Filecoin"
Set-CellText $ws "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText $ws "D33" "5.13"
Set-CellText $ws "E33" "  -2.96%  "
Set-CellText $ws "B34" "RenderToken"
Set-CellText $ws "C34" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws "D34" "4.66"
Set-CellText $ws "E34" "  +7.10%  "
Set-CellText $ws "B35" "WEMIXToken"
Set-CellText $ws "C35" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText $ws "D35" "2.48"
Set-CellText $ws "E35" "  -0.93%  "
Set-CellText $ws "B36" "Hedera"
Set-CellText $ws "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText $ws "D36" "0.0728"
Set-CellText $ws "E36" "  -2.07%  "
Set-CellText $ws "E37" "  -1.29%  "
Set-CellText $ws "B38" "LidoDAOToken"
Set-CellText $ws "C38" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-CellText $ws "D38" "2.95"
Set-CellText $ws "E38" "  -3.85%  "
Set-CellText $ws "B39" "ARBITRUM"
Set-CellText $ws "C39" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText $ws "D39" "1.86"
Set-CellText $ws "E39" "  -1.28%  "
Set-CellText $ws "B40" "Kaspa"
Set-CellText $ws "C40" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText $ws "D40" "0.103"
Set-CellText $ws "E40" "  -3.05%  "
Set-CellText $ws "B41" "Stellar"
Set-CellText $ws "C41" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText $ws "D41" "0.113"
Set-CellText $ws "E41" "  -2.40%  "
Set-CellText $ws "B42" "ApeXProtocol"
Set-CellText $ws "C42" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-CellText $ws "D42" "2.44"
Set-CellText $ws "E42" "  +5.48%  "
Set-CellText $ws "B43" "Maker"
Set-CellText $ws "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-CellText $ws "D43" "2.015.75"
Set-CellText $ws "E43" "  +2.25%  "
Set-CellText $ws "B44" "VeChain"
Set-CellText $ws "C44" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D44" "0.0286"
Set-CellText $ws "E44" "  -1.22%  "
Set-CellText $ws "B45" "EnergySwap"
Set-CellText $ws "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws "D45" "18.68"
Set-CellText $ws "E45" "  -2.80%  "
Set-CellText $ws "B46" "FraxShare"
Set-CellText $ws "C46" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText $ws "D46" "10.28"
Set-CellText $ws "E46" "  +3.62%  "
Set-CellText $ws "B47" "NEARProtocol"
Set-CellText $ws "C47" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-CellText $ws "D47" "2.92"
Set-CellText $ws "E47" "  -1.75%  "
Set-CellText $ws "B48" "MultiversX"
Set-CellText $ws "C48" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-CellText $ws "D48" "55.76"
Set-CellText $ws "E48" "  +0.39%  "
Set-CellText $ws "B49" "HuobiToken"
Set-CellText $ws "C49" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-CellText $ws "D49" "2.90"
Set-CellText $ws "E49" "  -0.29%  "
Set-CellText $ws "B50" "RocketPoolETH"
Set-CellText $ws "C50" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-CellText $ws "D50" "2.563.46"
Set-CellText $ws "E50" "  +1.06%  "
Set-CellText $ws "B51" "THORChain"
Set-CellText $ws "C51" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-CellText $ws "D51" "4.66"
Set-CellText $ws "E51" "  +1.58%  "
